$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the new date cells (A275:A288) with the same style (short-date format) as
# the existing date column by copying format from the row immediately above the
# new block; the actual date values are overwritten per-row below.
$ws.Range("A274").Copy($ws.Range("A275:A288"))

# Row 275
$ws.Range("A275").Value = 45826
$ws.Range("B275").Value = "Flowering"
$ws.Range("C275").Value = "Large"
$ws.Range("D275").Value = 69
$ws.Range("E275").Value = 82
$ws.Range("F275").Formula = "=ABS(D275-E275)"
$ws.Range("G275").Value = 3.75
$ws.Range("H275").Value = 0.6
$ws.Range("I275").Value = "No"
$ws.Range("J275").Value = 2
$ws.Range("K275").Value = "Dark"
$ws.Range("L275").Value = 8
$ws.Range("M275").Value = 0.82
$ws.Range("N275").Value = 66
$ws.Range("O275").Value = 29.85
$ws.Range("P275").Value = 12
$ws.Range("Q275").Value = 0.85
$ws.Range("R275").Value = 8.8000000000000007
$ws.Range("S275").Value = 35
$ws.Range("T275").Value = 24

# Row 276
$ws.Range("A276").Value = 45826
$ws.Range("B276").Value = "Nonflowering"
$ws.Range("C276").Value = "Medium"
$ws.Range("D276").Value = 69
$ws.Range("E276").Value = 82
$ws.Range("F276").Formula = "=ABS(D276-E276)"
$ws.Range("G276").Value = 3.75
$ws.Range("H276").Value = 0.9
$ws.Range("I276").Value = "No"
$ws.Range("J276").Value = 3
$ws.Range("K276").Value = "Dark"
$ws.Range("L276").Value = 8
$ws.Range("M276").Value = 0.82
$ws.Range("N276").Value = 66
$ws.Range("O276").Value = 29.85
$ws.Range("P276").Value = 12
$ws.Range("Q276").Value = 0.85
$ws.Range("R276").Value = 8.8000000000000007
$ws.Range("S276").Value = 35
$ws.Range("T276").Value = 24

# Row 277
$ws.Range("A277").Value = 45826
$ws.Range("B277").Value = "Nonflowering"
$ws.Range("C277").Value = "Small"
$ws.Range("D277").Value = 69
$ws.Range("E277").Value = 82
$ws.Range("F277").Formula = "=ABS(D277-E277)"
$ws.Range("G277").Value = 3.75
$ws.Range("H277").Value = 1
$ws.Range("I277").Value = "No"
$ws.Range("J277").Value = 3
$ws.Range("K277").Value = "Dark"
$ws.Range("L277").Value = 8
$ws.Range("M277").Value = 0.82
$ws.Range("N277").Value = 66
$ws.Range("O277").Value = 29.85
$ws.Range("P277").Value = 12
$ws.Range("Q277").Value = 0.85
$ws.Range("R277").Value = 8.8000000000000007
$ws.Range("S277").Value = 35
$ws.Range("T277").Value = 24

# Row 278
$ws.Range("A278").Value = 45826
$ws.Range("B278").Value = "Nonflowering"
$ws.Range("C278").Value = "Medium"
$ws.Range("D278").Value = 69
$ws.Range("E278").Value = 82
$ws.Range("F278").Formula = "=ABS(D278-E278)"
$ws.Range("G278").Value = 3.75
$ws.Range("H278").Value = 1.5
$ws.Range("I278").Value = "No"
$ws.Range("J278").Value = 3
$ws.Range("K278").Value = "Dark"
$ws.Range("L278").Value = 8
$ws.Range("M278").Value = 0.82
$ws.Range("N278").Value = 66
$ws.Range("O278").Value = 29.85
$ws.Range("P278").Value = 12
$ws.Range("Q278").Value = 0.85
$ws.Range("R278").Value = 8.8000000000000007
$ws.Range("S278").Value = 35
$ws.Range("T278").Value = 24

# Row 279
$ws.Range("A279").Value = 45826
$ws.Range("B279").Value = "Nonflowering"
$ws.Range("C279").Value = "Medium"
$ws.Range("D279").Value = 69
$ws.Range("E279").Value = 82
$ws.Range("F279").Formula = "=ABS(D279-E279)"
$ws.Range("G279").Value = 3.75
$ws.Range("H279").Formula = "=4/3"
$ws.Range("I279").Value = "No"
$ws.Range("J279").Value = 3
$ws.Range("K279").Value = "Dark"
$ws.Range("L279").Value = 8
$ws.Range("M279").Value = 0.82
$ws.Range("N279").Value = 66
$ws.Range("O279").Value = 29.85
$ws.Range("P279").Value = 12
$ws.Range("Q279").Value = 0.85
$ws.Range("R279").Value = 8.8000000000000007
$ws.Range("S279").Value = 35
$ws.Range("T279").Value = 24

# Row 280
$ws.Range("A280").Value = 45826
$ws.Range("B280").Value = "Nonflowering"
$ws.Range("C280").Value = "Large"
$ws.Range("D280").Value = 69
$ws.Range("E280").Value = 82
$ws.Range("F280").Formula = "=ABS(D280-E280)"
$ws.Range("G280").Value = 3.75
$ws.Range("H280").Value = 4
$ws.Range("I280").Value = "No"
$ws.Range("J280").Value = 4
$ws.Range("K280").Value = "Dark"
$ws.Range("L280").Value = 8
$ws.Range("M280").Value = 0.82
$ws.Range("N280").Value = 66
$ws.Range("O280").Value = 29.85
$ws.Range("P280").Value = 12
$ws.Range("Q280").Value = 0.85
$ws.Range("R280").Value = 8.8000000000000007
$ws.Range("S280").Value = 35
$ws.Range("T280").Value = 24

# Row 281
$ws.Range("A281").Value = 45826
$ws.Range("B281").Value = "Tree"
$ws.Range("C281").Value = "Medium"
$ws.Range("D281").Value = 69
$ws.Range("E281").Value = 82
$ws.Range("F281").Formula = "=ABS(D281-E281)"
$ws.Range("G281").Value = 3.75
$ws.Range("H281").Value = 5.6
$ws.Range("I281").Value = "No"
$ws.Range("J281").Value = 1
$ws.Range("K281").Value = "Dark"
$ws.Range("L281").Value = 8
$ws.Range("M281").Value = 0.82
$ws.Range("N281").Value = 66
$ws.Range("O281").Value = 29.85
$ws.Range("P281").Value = 12
$ws.Range("Q281").Value = 0.85
$ws.Range("R281").Value = 8.8000000000000007
$ws.Range("S281").Value = 35
$ws.Range("T281").Value = 24

# Row 282
$ws.Range("A282").Value = 45827
$ws.Range("B282").Value = "Flowering"
$ws.Range("C282").Value = "Large"
$ws.Range("D282").Value = 58
$ws.Range("E282").Value = 78
$ws.Range("F282").Formula = "=ABS(D282-E282)"
$ws.Range("G282").Value = 1.1100000000000001
$ws.Range("H282").Value = 0.1
$ws.Range("I282").Value = "No"
$ws.Range("J282").Value = 2
$ws.Range("K282").Value = "Dark"
$ws.Range("L282").Value = 8
$ws.Range("M282").Value = 0.82
$ws.Range("N282").Value = 66
$ws.Range("O282").Value = 29.76
$ws.Range("P282").Value = 27
$ws.Range("Q282").Value = 0.9
$ws.Range("R282").Value = 8.6999999999999993
$ws.Range("S282").Value = 37
$ws.Range("T282").Value = 13

# Row 283
$ws.Range("A283").Value = 45827
$ws.Range("B283").Value = "Nonflowering"
$ws.Range("C283").Value = "Medium"
$ws.Range("D283").Value = 58
$ws.Range("E283").Value = 78
$ws.Range("F283").Formula = "=ABS(D283-E283)"
$ws.Range("G283").Value = 1.1100000000000001
$ws.Range("H283").Value = 0.25
$ws.Range("I283").Value = "No"
$ws.Range("J283").Value = 3
$ws.Range("K283").Value = "Dark"
$ws.Range("L283").Value = 8
$ws.Range("M283").Value = 0.82
$ws.Range("N283").Value = 66
$ws.Range("O283").Value = 29.76
$ws.Range("P283").Value = 27
$ws.Range("Q283").Value = 0.9
$ws.Range("R283").Value = 8.6999999999999993
$ws.Range("S283").Value = 37
$ws.Range("T283").Value = 13

# Row 284
$ws.Range("A284").Value = 45827
$ws.Range("B284").Value = "Nonflowering"
$ws.Range("C284").Value = "Small"
$ws.Range("D284").Value = 58
$ws.Range("E284").Value = 78
$ws.Range("F284").Formula = "=ABS(D284-E284)"
$ws.Range("G284").Value = 1.1100000000000001
$ws.Range("H284").Value = 0.4
$ws.Range("I284").Value = "No"
$ws.Range("J284").Value = 3
$ws.Range("K284").Value = "Dark"
$ws.Range("L284").Value = 8
$ws.Range("M284").Value = 0.82
$ws.Range("N284").Value = 66
$ws.Range("O284").Value = 29.76
$ws.Range("P284").Value = 27
$ws.Range("Q284").Value = 0.9
$ws.Range("R284").Value = 8.6999999999999993
$ws.Range("S284").Value = 37
$ws.Range("T284").Value = 13

# Row 285
$ws.Range("A285").Value = 45827
$ws.Range("B285").Value = "Nonflowering"
$ws.Range("C285").Value = "Medium"
$ws.Range("D285").Value = 58
$ws.Range("E285").Value = 78
$ws.Range("F285").Formula = "=ABS(D285-E285)"
$ws.Range("G285").Value = 1.1100000000000001
$ws.Range("H285").Value = 0.65
$ws.Range("I285").Value = "No"
$ws.Range("J285").Value = 3
$ws.Range("K285").Value = "Dark"
$ws.Range("L285").Value = 8
$ws.Range("M285").Value = 0.82
$ws.Range("N285").Value = 66
$ws.Range("O285").Value = 29.76
$ws.Range("P285").Value = 27
$ws.Range("Q285").Value = 0.9
$ws.Range("R285").Value = 8.6999999999999993
$ws.Range("S285").Value = 37
$ws.Range("T285").Value = 13

# Row 286
$ws.Range("A286").Value = 45827
$ws.Range("B286").Value = "Nonflowering"
$ws.Range("C286").Value = "Medium"
$ws.Range("D286").Value = 58
$ws.Range("E286").Value = 78
$ws.Range("F286").Formula = "=ABS(D286-E286)"
$ws.Range("G286").Value = 1.1100000000000001
$ws.Range("H286").Value = 0.7
$ws.Range("I286").Value = "No"
$ws.Range("J286").Value = 3
$ws.Range("K286").Value = "Dark"
$ws.Range("L286").Value = 8
$ws.Range("M286").Value = 0.82
$ws.Range("N286").Value = 66
$ws.Range("O286").Value = 29.76
$ws.Range("P286").Value = 27
$ws.Range("Q286").Value = 0.9
$ws.Range("R286").Value = 8.6999999999999993
$ws.Range("S286").Value = 37
$ws.Range("T286").Value = 13

# Row 287
$ws.Range("A287").Value = 45827
$ws.Range("B287").Value = "Nonflowering"
$ws.Range("C287").Value = "Large"
$ws.Range("D287").Value = 58
$ws.Range("E287").Value = 78
$ws.Range("F287").Formula = "=ABS(D287-E287)"
$ws.Range("G287").Value = 1.1100000000000001
$ws.Range("H287").Value = 1.5
$ws.Range("I287").Value = "No"
$ws.Range("J287").Value = 4
$ws.Range("K287").Value = "Dark"
$ws.Range("L287").Value = 8
$ws.Range("M287").Value = 0.82
$ws.Range("N287").Value = 66
$ws.Range("O287").Value = 29.76
$ws.Range("P287").Value = 27
$ws.Range("Q287").Value = 0.9
$ws.Range("R287").Value = 8.6999999999999993
$ws.Range("S287").Value = 37
$ws.Range("T287").Value = 13

# Row 288
$ws.Range("A288").Value = 45827
$ws.Range("B288").Value = "Tree"
$ws.Range("C288").Value = "Medium"
$ws.Range("D288").Value = 58
$ws.Range("E288").Value = 78
$ws.Range("F288").Formula = "=ABS(D288-E288)"
$ws.Range("G288").Value = 1.1100000000000001
$ws.Range("H288").Value = 3
$ws.Range("I288").Value = "No"
$ws.Range("J288").Value = 1
$ws.Range("K288").Value = "Dark"
$ws.Range("L288").Value = 8
$ws.Range("M288").Value = 0.82
$ws.Range("N288").Value = 66
$ws.Range("O288").Value = 29.76
$ws.Range("P288").Value = 27
$ws.Range("Q288").Value = 0.9
$ws.Range("R288").Value = 8.6999999999999993
$ws.Range("S288").Value = 37
$ws.Range("T288").Value = 13

# Match the saved workbook's selection cursor recorded in the edit.
[void]$ws.Range("U3").Select()
